# B6-PowerPoint.pptx edit
#
# 1) Three tables (on what end up being slides 14, 15 and 16) get their
#    table style switched from the custom "Table_0" style
#    ({3615A551-651C-4BA5-91F9-9C423CB14D86}) to a different (built-in)
#    table style ({C75D9288-0F8C-4855-9AB3-EC7FA5896237}).
#
# 2) The deck's theme ends up using the plain "Office" colour palette
#    (it currently uses the "Integral"/"Red Violet" palette) - i.e. the
#    Design/Theme applied to the slide master (and therefore to every
#    slide, since the colours are all inherited via schemeClr) is
#    switched back to the default Office colours.

$p = $ppt.ActivePresentation

# --- 1. Retarget the three affected tables' style -----------------------

$targetStyleId = "{C75D9288-0F8C-4855-9AB3-EC7FA5896237}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2. Swap the active theme's colour scheme back to "Office" ----------

$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# RGB() packs as 0xBBGGRR for the COM long, matching VBA's RGB(r,g,b).
function RGB($r, $g, $b) { return $b * 65536 + $g * 256 + $r }

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - the standard "Office"
# theme colour scheme (as opposed to the current "Red Violet" / Integral
# scheme).
$colors.Item(1).RGB  = RGB 0x00 0x00 0x00   # tx1  / dk1      - 000000
$colors.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # bg1  / lt1      - FFFFFF
$colors.Item(3).RGB  = RGB 0x44 0x54 0x6A   # tx2  / dk2      - 44546A
$colors.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # bg2  / lt2      - E7E6E6
$colors.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1         - 5B9BD5
$colors.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2         - ED7D31
$colors.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3         - A5A5A5
$colors.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4         - FFC000
$colors.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5         - 4472C4
$colors.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6         - 70AD47
$colors.Item(11).RGB = RGB 0x05 0x63 0xC1   # hyperlink       - 0563C1
$colors.Item(12).RGB = RGB 0x95 0x4F 0x72   # followed hlink  - 954F72

# Best-effort: try to relabel the theme/colour-scheme names too (no-ops on
# hosts that treat Theme/Design names as read-only, harmless either way).
$design = $p.Designs.Item(1)
$design.Name = "Office Theme"
$master.Theme.Name = "Office Theme"
$colors.Name = "Office"
$p.TemplateName = "Office Theme"
